# Update students_data.xlsx:
#  - fix two existing "Fees/Admission" date strings (registration #1 and #2)
#  - fill in the previously-empty Address for registration #2
#  - set a proper "Fees Paid Till Date" for registration #2
#  - append a brand-new student row: registration #3 ("ashish")
#
# This sheet stores every column (including numbers/dates/phone numbers) as
# plain text, not native Excel numbers/dates, so numeric-looking values are
# written with a leading apostrophe (the same trick Excel's UI uses) to keep
# them as literal text instead of being auto-converted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# --- Fix existing values -------------------------------------------------

# Registration #1 ("satyam") - Fees Paid Till Date (M2)
$ws.Range("M2").Value = "'2025-03-03"

# Registration #2 ("test") - Admission Date (B3)
$ws.Range("B3").Value = "'2025-02-03"

# Registration #2 ("test") - Address was missing, now filled in (E3)
$ws.Range("E3").Value = "testaddress"

# Registration #2 ("test") - Locker Number changed from 3 to 2 (L3)
$ws.Range("L3").Value = "'2"

# Registration #2 ("test") - Fees Paid Till Date (M3)
$ws.Range("M3").Value = "'2025-03-13"

# --- Append new student row (row 4) --------------------------------------

$ws.Range("A4").Value = "'3"
$ws.Range("B4").Value = "'2025-01-01"
$ws.Range("C4").Value = "ashish"
$ws.Range("D4").Value = "ashish ke babuji"
$ws.Range("E4").Value = "ramkrishnanagar"
$ws.Range("F4").Value = "'7250585059"
$ws.Range("G4").Value = "06:00-10:00"
$ws.Range("H4").Value = "'1"
$ws.Range("I4").Value = "'0"
$ws.Range("J4").Value = "'350.00"
$ws.Range("K4").Value = "'150.00"
$ws.Range("L4").Value = "'6"
$ws.Range("M4").Value = "'2025-02-14"
